$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "HDJ" labels to have a space between the word and the number
$ws.Range("A15").Value = "HDJ 1"
$ws.Range("A16").Value = "HDJ 2"
$ws.Range("A17").Value = "HDJ 3"

# Update the selected range on the sheet (as saved in the workbook view state)
$ws.Range("A18").Select()
